# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2-143) from 45181 (2023-09-12) to 45182 (2023-09-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C143").Value2 = 45182
